$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.656.39"
$ws.Range("E2").Value = "  -2.06%  "
$ws.Range("D3").Value = "3.744.18"
$ws.Range("E3").Value = "  +2.01%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'619.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.12%  "
$ws.Range("D6").Value = "'179.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.67%  "
$ws.Range("D7").Value = "3.738.03"
$ws.Range("E7").Value = "  +1.96%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").Value = "'0.530"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.62%  "
$ws.Range("D10").Value = "'0.168"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.93%  "
$ws.Range("D11").Value = "'6.30"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.15%  "
$ws.Range("D12").Value = "'0.489"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.15%  "
$ws.Range("D13").Value = "'41.00"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.00%  "
$ws.Range("D14").Value = "'0.0000259"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.53%  "
$ws.Range("D15").Value = "4.375.18"
$ws.Range("E15").Value = "  +2.22%  "
$ws.Range("D16").Value = "3.754.99"
$ws.Range("E16").Value = "  +2.26%  "
$ws.Range("D17").Value = "69.814.59"
$ws.Range("E17").Value = "  -1.80%  "
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("D19").Value = "'7.61"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.21%  "
$ws.Range("D20").Value = "'16.70"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.38%  "
$ws.Range("D21").Value = "'506.13"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.21%  "
$ws.Range("D22").Value = "'9.52"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.34%  "
$ws.Range("D23").Value = "'0.727"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.33%  "
$ws.Range("D24").Value = "'2.50"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.23%  "
$ws.Range("D25").Value = "'87.12"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.60%  "
$ws.Range("D26").Value = "'13.13"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.39%  "
$ws.Range("D27").Value = "'11.13"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.72%  "
$ws.Range("E28").Value = "  +24.05%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").Value = "'2.51"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.66%  "
$ws.Range("D31").Value = "'2.89"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.15%  "
$ws.Range("D32").Value = "'7.90"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.70%  "
$ws.Range("D33").Value = "'31.01"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.38%  "
$ws.Range("E34").Value = "  -0.76%  "
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("E36").Value = "  +4.19%  "
$ws.Range("D37").Value = "'6.20"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.87%  "
$ws.Range("D38").Value = "'0.337"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.48%  "
$ws.Range("E39").Value = "  +2.17%  "
$ws.Range("D40").Value = "'2.11"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.00%  "
$ws.Range("D41").Value = "'50.16"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.37%  "
$ws.Range("D42").Value = "'45.67"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.84%  "
$ws.Range("D43").Value = "'428.14"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.79%  "
$ws.Range("D44").Value = "'8.72"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.26%  "
$ws.Range("D45").Value = "'2.86"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.06%  "
$ws.Range("D46").Value = "3.000.90"
$ws.Range("E46").Value = "  -3.97%  "
$ws.Range("D47").Value = "'0.0364"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.07%  "
$ws.Range("D48").Value = "'27.36"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.73%  "
$ws.Range("E49").Value = "  -0.06%  "
$ws.Range("D50").Value = "'136.91"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.61%  "
$ws.Range("D51").Value = "'2.49"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.12%  "
